$wb = $excel.ActiveWorkbook

# --- 1. Update the status text from "Ready for handoff" to "In Translation" ---
# This shared string is referenced from:
#   Overview!E2, Overview!F2, zh-cn!C2, de-de!C2
# Setting the value on each of these cells updates every occurrence.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2. Shrink the now-narrower status columns to fit the shorter text ---
# Overview columns E (zh-cn status) and F (de-de status)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
